$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1888
$ws1.Range("F3").Value = 500
$ws1.Range("F6").Value = 2639
$ws1.Range("F7").Value = 179
$ws1.Range("F8").Value = 94
$ws1.Range("F10").Value = 1556
$ws1.Range("F11").Value = 540
$ws1.Range("F21").Value = 192
$ws1.Range("F23").Value = 1693
$ws1.Range("F26").Value = 25
$ws1.Range("F28").Value = 213
$ws1.Range("F30").Value = 431

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1888
$ws4.Range("F4").Value = 500
$ws4.Range("F7").Value = 2639
$ws4.Range("F8").Value = 179
$ws4.Range("F9").Value = 94
$ws4.Range("F11").Value = 1556
$ws4.Range("F12").Value = 540
$ws4.Range("F22").Value = 192
$ws4.Range("F24").Value = 1693
$ws4.Range("F27").Value = 25
$ws4.Range("F29").Value = 213
$ws4.Range("F31").Value = 431
